$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its text representation (values look numeric but
# must remain stored as text, matching the original inline-string cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.927.13"
$ws.Range("E2").Value = "  +6.32%  "
$ws.Range("D3").Value = "1.877.38"
$ws.Range("E3").Value = "  +5.44%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("D5").Value = "248.42"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.4998"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("D8").Value = "45.84"
$ws.Range("E8").Value = "  +8.11%  "
$ws.Range("D9").Value = "0.2854"
$ws.Range("E9").Value = "  +6.73%  "
$ws.Range("D10").Value = "0.06539"
$ws.Range("E10").Value = "  +4.30%  "
$ws.Range("D11").Value = "1.880.62"
$ws.Range("E11").Value = "  +5.64%  "
$ws.Range("D12").Value = "17.08"
$ws.Range("E12").Value = "  +3.75%  "
$ws.Range("D13").Value = "0.07240"
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("D14").Value = "0.6673"
$ws.Range("E14").Value = "  +6.35%  "
$ws.Range("D15").Value = "85.02"
$ws.Range("E15").Value = "  +6.37%  "
$ws.Range("D16").Value = "4.829"
$ws.Range("E16").Value = "  +3.69%  "
$ws.Range("D17").Value = "29.939.25"
$ws.Range("E17").Value = "  +6.45%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "12.90"
$ws.Range("E19").Value = "  +7.11%  "
$ws.Range("D20").Value = "0.000007512"
$ws.Range("E20").Value = "  +3.59%  "
$ws.Range("D21").Value = "0.9994"
$ws.Range("D22").Value = "2.122.88"
$ws.Range("E22").Value = "  +5.68%  "
$ws.Range("D23").Value = "4.765"
$ws.Range("E23").Value = "  +4.40%  "
$ws.Range("E24").Value = "  +5.19%  "
$ws.Range("D25").Value = "9.011"
$ws.Range("E25").Value = "  +3.22%  "
$ws.Range("D26").Value = "145.25"
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("E27").Value = "  +23.45%  "
$ws.Range("D28").Value = "16.71"
$ws.Range("E28").Value = "  +6.13%  "
$ws.Range("D29").Value = "1.952"
$ws.Range("E29").Value = "  +5.02%  "
$ws.Range("D30").Value = "1.371"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("D31").Value = "4.197"
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("D32").Value = "0.08657"
$ws.Range("E32").Value = "  +4.76%  "
$ws.Range("D33").Value = "3.891"
$ws.Range("E33").Value = "  +3.61%  "
$ws.Range("D34").Value = "0.05064"
$ws.Range("E34").Value = "  +3.62%  "
$ws.Range("D35").Value = "1.131"
$ws.Range("E35").Value = "  +5.35%  "
$ws.Range("D36").Value = "0.6888"
$ws.Range("E36").Value = "  +5.85%  "
$ws.Range("D37").Value = "2.683"
$ws.Range("E37").Value = "  +2.38%  "
$ws.Range("D38").Value = "2.294"
$ws.Range("E38").Value = "  +12.38%  "
$ws.Range("E39").Value = "  +5.56%  "
$ws.Range("D40").Value = "0.9605"
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("D41").Value = "0.01634"
$ws.Range("E41").Value = "  +5.44%  "
$ws.Range("E42").Value = "  +3.56%  "
$ws.Range("D43").Value = "104.52"
$ws.Range("E43").Value = "  +4.58%  "
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "0.4214"
$ws.Range("E45").Value = "  +5.75%  "
$ws.Range("D46").Value = "7.447"
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("D47").Value = "0.1254"
$ws.Range("E47").Value = "  +3.36%  "
$ws.Range("D48").Value = "0.05633"
$ws.Range("E48").Value = "  +3.65%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.291"
$ws.Range("E49").Value = "  +3.41%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "32.39"
$ws.Range("E50").Value = "  +5.33%  "
$ws.Range("D51").Value = "0.3710"
$ws.Range("E51").Value = "  +6.63%  "

# Restore default (unstyled) cell style now that the text values are set,
# so the text-number-format style does not linger on the cells.
$ws.Range("D2:D51").Style = "Normal"

